$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Widen column I (Status) from 10 to 14 (raw OOXML width units).
#    COM ColumnWidth is offset by -0.83 from the stored OOXML width for this font.
$ws.Columns.Item(9).ColumnWidth = 13.17

# 2. Swap the "Recorded By" ordering everywhere it appears:
#    "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
$ws.Cells.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System")

# 3. The 12 "SURGERY SEMINAR/SLIDE" sessions (27/12/2025 08:30) that were
#    still shown as "Pending" have now been marked "Not Recorded".
$notRecordedRows = 16, 35, 54, 73, 91, 109, 127, 145, 163, 182, 201, 220
foreach ($r in $notRecordedRows) {
    $ws.Range("I$r").Value = "Not Recorded"
}

# 4. Roll the 12 newly-flagged sessions from "Pending" into "Missing" in the
#    overall Class Statistics block.
$ws.Range("L7").Value = 12
$ws.Range("L8").Value = 132

# 5. Mirror that same move (one session each) in the per-group Group
#    Statistics table for every B1 group (rows 15-26: Missing +1, Pending -1).
for ($r = 15; $r -le 26; $r++) {
    $ws.Range("P$r").Value = 1
    $ws.Range("Q$r").Value = 11
}
